$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header for new column D
$ws.Range("D2").Value = "color"

# Set "zinc" for rows corresponding to IS1-77-001 (row 6) and SFTMH-23T (row 14)
$ws.Range("D6").Value = "zinc"
$ws.Range("D14").Value = "zinc"

# Apply bordered style to D2:D18 to match the new cellXfs style (borderId=1)
for ($r = 2; $r -le 18; $r++) {
  $cell = $ws.Range("D$r")
  $cell.Borders.Color = 0
  $cell.Borders.LineStyle = 1
  $cell.Borders.Weight = 2
}

# Selection matches the final state in the diff
$ws.Range("D15").Select()
